$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking price values that must remain as text
$textCells = @("D5","D6","D10","D11","D12","D16","D19","D20","D21","D22","D23","D24","D25","D26","D27","D29","D32","D33","D34","D35","D36","D41","D42","D43","D44","D46","D47","D48","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '62.933.38'
$ws.Range('E2').Value = '  -5.34%  '
$ws.Range('D3').Value = '3.117.89'
$ws.Range('E3').Value = '  -5.82%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').Value = '559.87'
$ws.Range('E5').Value = '  -4.71%  '
$ws.Range('D6').Value = '162.78'
$ws.Range('E6').Value = '  -9.95%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('E8').Value = '  -10.01%  '
$ws.Range('D9').Value = '3.116.44'
$ws.Range('E9').Value = '  -5.70%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').Value = '0.117'
$ws.Range('E10').Value = '  -8.07%  '
$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').Value = '6.70'
$ws.Range('E11').Value = '  -2.12%  '
$ws.Range('D12').Value = '0.380'
$ws.Range('E12').Value = '  -5.95%  '
$ws.Range('D13').Value = '3.664.97'
$ws.Range('E13').Value = '  -5.69%  '
$ws.Range('E14').Value = '  -2.28%  '
$ws.Range('D15').Value = '63.113.93'
$ws.Range('E15').Value = '  -5.02%  '
$ws.Range('D16').Value = '24.64'
$ws.Range('E16').Value = '  -7.82%  '
$ws.Range('D17').Value = '3.128.50'
$ws.Range('E17').Value = '  -5.07%  '
$ws.Range('E18').Value = '  -5.93%  '
$ws.Range('D19').Value = '406.58'
$ws.Range('E19').Value = '  -4.04%  '
$ws.Range('D20').Value = '12.60'
$ws.Range('E20').Value = '  -4.08%  '
$ws.Range('D21').Value = '5.17'
$ws.Range('E21').Value = '  -5.53%  '
$ws.Range('D22').Value = '7.05'
$ws.Range('E22').Value = '  -3.90%  '
$ws.Range('D23').Value = '1.00'
$ws.Range('E23').Value = '  +0.09%  '
$ws.Range('D24').Value = '5.69'
$ws.Range('E24').Value = '  +0.37%  '
$ws.Range('D25').Value = '68.60'
$ws.Range('E25').Value = '  -4.04%  '
$ws.Range('D26').Value = '0.200'
$ws.Range('E26').Value = '  -2.60%  '
$ws.Range('D27').Value = '0.490'
$ws.Range('E27').Value = '  -4.82%  '
$ws.Range('E28').Value = '  -11.07%  '
$ws.Range('D29').Value = '8.72'
$ws.Range('E29').Value = '  -4.30%  '
$ws.Range('E30').Value = '  +0.14%  '
$ws.Range('E31').Value = '  -0.09%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '21.33'
$ws.Range('E32').Value = '  -4.84%  '
$ws.Range('B33').Value = 'PancakeSwap'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D33').Value = '1.79'
$ws.Range('E33').Value = '  -7.08%  '
$ws.Range('D34').Value = '4.90'
$ws.Range('E34').Value = '  -5.17%  '
$ws.Range('D35').Value = '6.26'
$ws.Range('E35').Value = '  -5.29%  '
$ws.Range('D36').Value = '154.37'
$ws.Range('E36').Value = '  -3.63%  '
$ws.Range('E37').Value = '  -6.92%  '
$ws.Range('E38').Value = '  -6.99%  '
$ws.Range('D39').Value = '2.718.28'
$ws.Range('E39').Value = '  -4.89%  '
$ws.Range('E40').Value = '  -8.46%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = '4.13'
$ws.Range('E41').Value = '  -4.71%  '
$ws.Range('B42').Value = 'EnergySwap'
$ws.Range('C42').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D42').Value = '23.73'
$ws.Range('E42').Value = '  -9.89%  '
$ws.Range('D43').Value = '38.30'
$ws.Range('E43').Value = '  -3.54%  '
$ws.Range('D44').Value = '0.695'
$ws.Range('E44').Value = '  -7.57%  '
$ws.Range('E45').Value = '  -7.93%  '
$ws.Range('D46').Value = '0.0257'
$ws.Range('E46').Value = '  -5.78%  '
$ws.Range('D47').Value = '5.22'
$ws.Range('E47').Value = '  -11.61%  '
$ws.Range('B48').Value = 'Bittensor'
$ws.Range('C48').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D48').Value = '286.39'
$ws.Range('E48').Value = '  -8.18%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').Value = '20.99'
$ws.Range('E49').Value = '  -9.11%  '
$ws.Range('E50').Value = '  +0.00%  '
$ws.Range('D51').Value = '0.0973'
$ws.Range('E51').Value = '  -6.39%  '

# restore default style (remove explicit number-format style) on affected cells
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}